$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 block - labels reuse shared strings 0,1,2 whose text changes
$ws.Range("E7").Value = "18-23"
$ws.Range("F7").Value = "14-17"
$ws.Range("G7").Value = "10-13"
$ws.Range("H7").Value = 45541

# Row 10 block
$ws.Range("E10").Value = "14-23"
$ws.Range("G10").Value = "10-13"
$ws.Range("H10").Value = 45541

# Row 13 block
$ws.Range("E13").Value = "14-23"
$ws.Range("G13").Value = "10-13"
$ws.Range("H13").Value = 45572
$ws.Range("I13").Value = 45446

# Row 14
$ws.Range("D14").Value = "S"
$ws.Range("K14").Value = "SW"

# Row 16 block
$ws.Range("E16").Value = "14-23"
$ws.Range("G16").Value = "10-13"
$ws.Range("H16").Value = 45541

# Row 17
$ws.Range("D17").Value = "B"
$ws.Range("K17").Value = "BEQ, BNE, BLT, BGE"

# Row 19
$ws.Range("E19").Value = "7-23"
$ws.Range("I19").Value = 45446

# Row 20
$ws.Range("D20").Value = "J"
$ws.Range("K20").Value = "JAL"
